# Updated cryptos list on Tue Dec 26 22:43:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

Set-TextCell "D2" "42.663.19"
Set-TextCell "E2" "  -2.69%  "
Set-TextCell "D3" "2.227.46"
Set-TextCell "E3" "  -2.54%  "
Set-TextCell "E4" "  +0.15%  "
Set-TextCell "D5" "112.28"
Set-TextCell "E5" "  -6.80%  "
Set-TextCell "D6" "297.21"
Set-TextCell "D7" "0.621"
Set-TextCell "E7" "  -4.05%  "
Set-TextCell "E8" "  -0.39%  "
Set-TextCell "D9" "0.609"
Set-TextCell "E9" "  -4.24%  "
Set-TextCell "D10" "44.35"
Set-TextCell "E10" "  -8.64%  "
Set-TextCell "E11" "  -3.65%  "
Set-TextCell "D12" "54.47"
Set-TextCell "E12" "  +0.11%  "
Set-TextCell "D13" "8.89"
Set-TextCell "E13" "  -4.49%  "
Set-TextCell "D14" "1.01"
Set-TextCell "E14" "  +8.97%  "
Set-TextCell "D15" "0.103"
Set-TextCell "E15" "  -3.12%  "
Set-TextCell "D16" "15.13"
Set-TextCell "E16" "  -3.09%  "
Set-TextCell "D17" "2.559.99"
Set-TextCell "E17" "  -2.66%  "
Set-TextCell "D18" "2.233.45"
Set-TextCell "E18" "  -1.72%  "
Set-TextCell "D19" "42.473.46"
Set-TextCell "E19" "  -3.04%  "
Set-TextCell "D20" "7.39"
Set-TextCell "E20" "  +5.85%  "
Set-TextCell "E21" "  -4.16%  "
Set-TextCell "D22" "72.89"
Set-TextCell "E22" "  +0.65%  "
Set-TextCell "D23" "3.51"
Set-TextCell "E23" "  +21.60%  "
Set-TextCell "E24" "  -2.61%  "
Set-TextCell "D25" "229.63"
Set-TextCell "E25" "  -2.93%  "
Set-TextCell "D26" "9.25"
Set-TextCell "E26" "  -4.34%  "
Set-TextCell "D27" "11.74"
Set-TextCell "E27" "  -2.80%  "
Set-TextCell "D28" "0.999"
Set-TextCell "E28" "  -1.63%  "
Set-TextCell "E29" "  -0.62%  "
Set-TextCell "D30" "38.51"
Set-TextCell "E30" "  -10.12%  "
Set-TextCell "D31" "3.25"
Set-TextCell "E31" "  -4.09%  "
Set-TextCell "D32" "173.82"
Set-TextCell "D33" "21.11"
Set-TextCell "E33" "  -2.85%  "
Set-TextCell "D34" "0.0905"
Set-TextCell "E34" "  -2.86%  "
Set-TextCell "D35" "5.20"
Set-TextCell "E35" "  +12.57%  "
Set-TextCell "E36" "  -2.56%  "
Set-TextCell "D37" "4.32"
Set-TextCell "E38" "  -4.00%  "
Set-TextCell "E39" "  -2.76%  "
Set-TextCell "E40" "  -2.81%  "
Set-TextCell "E41" "  -5.75%  "
Set-TextCell "D42" "72.16"
Set-TextCell "E42" "  -3.14%  "
Set-TextCell "E43" "  -1.15%  "
Set-TextCell "D44" "12.81"
Set-TextCell "E44" "  -7.21%  "
Set-TextCell "E45" "  +0.13%  "
Set-TextCell "D46" "1.33"
Set-TextCell "E46" "  -4.95%  "
Set-TextCell "E47" "  -6.47%  "
Set-TextCell "E48" "  +5.02%  "
Set-TextCell "D49" "103.46"
Set-TextCell "E49" "  +0.95%  "
Set-TextCell "D50" "8.56"
Set-TextCell "E50" "  -0.01%  "
Set-TextCell "E51" "  +6.14%  "
